$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear Source Well (C) and Plate (E) values for rows 4 and 5,
# leaving the cells blank but formatted the same.
$ws.Range("C4").Value = $null
$ws.Range("E4").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("E5").Value = $null

# Widen column D (target stored width 24.1640625; engine quantizes
# ColumnWidth to 1/6 steps, so feed the closest achievable input).
$ws.Columns.Item(4).ColumnWidth = 23.3333333333333

# Move the active cell selection to A5 instead of A6.
$ws.Range("A5").Select()
